# Word has a 40 character limit for bookmark names, and bookmark names
# must begin with a letter.  The two headings below use long,
# auto-generated identifiers that violate those rules, so rename the
# matching bookmarks (and the hyperlinks that point at them) to the
# SHA1-hash-based names pandoc now falls back to.
$d = $word.ActiveDocument

$renames = @{
    "remote-folder-or-longlonglonglonglong-file-with-manymanymanymany-letters-inside-opening" = "X49da2d776f7a640cd76098979e5788f8119bc44";
    "remote-folder-or-longlonglonglonglong-file-with-manymanymanymany-letters-inside-closing" = "Xb95b585046f38c7739779215f99b6b21152b861"
}

# --- Rename the bookmarks that mark the long headings ---
# Bookmark.Name isn't directly settable through this object model, so
# capture the bookmark's range, delete it, and re-add it under the new
# name.  Collect everything first since mutating Bookmarks while
# enumerating it is unsafe.
$bookmarksToFix = @()
foreach ($bm in $d.Bookmarks) {
    if ($renames.ContainsKey($bm.Name)) {
        $bookmarksToFix += ,@($bm.Name, $renames[$bm.Name], $bm.Start, $bm.End)
    }
}

foreach ($item in $bookmarksToFix) {
    $oldName = $item[0]
    $newName = $item[1]
    $rng = $d.Range($item[2], $item[3])
    $d.Bookmarks($oldName).Delete()
    $d.Bookmarks.Add($newName, $rng)
}

# --- Point every hyperlink that referenced the old anchors at the new ones ---
foreach ($hl in $d.Hyperlinks) {
    if ($renames.ContainsKey($hl.SubAddress)) {
        $hl.SubAddress = $renames[$hl.SubAddress]
    }
}
